$d = $word.ActiveDocument

# Locate the paragraph that currently holds the long "Ahora hay que
# modificar..." text (it ends the document body, right before sectPr).
$n = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs.Item($i).Range.Text
    if ($txt -like "Ahora hay que modificar las migraciones*") {
        $n = $i
        break
    }
}

if ($n -eq -1) {
    throw "Could not locate target paragraph"
}

# Clear the paragraph's text but keep its paragraph mark (and therefore its
# pPr / mark-formatting rPr) intact -- this becomes the first, now-empty,
# paragraph of the split.
$p = $d.Paragraphs.Item($n)
$r = $p.Range
$clearRange = $d.Range($r.Start, $r.End - 1)
$clearRange.Text = ""

# --- Second paragraph: "Hacer métodos de controladores." -------------------
$p = $d.Paragraphs.Item($n)
$p.Range.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item($n + 1)
$p2.Range.InsertAfter("Hacer métodos de controladores.")

# --- Third paragraph: "Hacer vistas." (keeps the _GoBack bookmark, placed
#     right after the run's text, before the paragraph mark) ---------------
$p2 = $d.Paragraphs.Item($n + 1)
$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item($n + 2)
# A trailing placeholder character is inserted temporarily: adding a
# collapsed bookmark exactly at "end of paragraph text" (i.e. immediately
# before the paragraph mark) is mishandled unless there is still a
# character after that position, so we add one, place the bookmark, then
# remove it again.
$p3.Range.InsertAfter("Hacer vistas.X")

$p3 = $d.Paragraphs.Item($n + 2)
$bmPos = $p3.Range.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bmRange)

$p3 = $d.Paragraphs.Item($n + 2)
$placeholder = $d.Range($p3.Range.End - 2, $p3.Range.End - 1)
$placeholder.Text = ""

# --- Fourth paragraph: "Hacer panel de administración..." ------------------
$p3 = $d.Paragraphs.Item($n + 2)
$p3.Range.InsertParagraphAfter()

$p4 = $d.Paragraphs.Item($n + 3)
$p4.Range.InsertAfter("Hacer panel de administración para usuario registrado y para administrador.")
